# Entities and Relationships.docx - apply the commit's changes via Word COM interop.
#
# Summary of the edit:
#  1. Remove the stray __DdeLink__ bookmark wrapped around the first
#     "Entities:" heading (bookmarkStart/bookmarkEnd only, text untouched).
#  2. After the "Route(...)" entity paragraph, add three more entity
#     paragraphs: Town(name, garage) / Garage(buses) / Stage(drivers).
#  3. Replace the two terse relationship lines ("Drivers drive Buses",
#     "Buses drive Routes") with the fuller relationship list, including
#     tab-separated cardinality annotations.

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Common run / paragraph-mark formatting used throughout this document body:
# Ubuntu font, non-bold, 18pt (sz 36 half-points). The pPr copy carries the
# (harmless) duplicated <w:b w:val="false"/> that the source document itself
# uses for every paragraph mark's rPr.
$rPr36 = '<w:rPr><w:rFonts w:ascii="Ubuntu" w:hAnsi="Ubuntu"/><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr>'
$pPr36 = '<w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:rFonts w:ascii="Ubuntu" w:hAnsi="Ubuntu"/><w:b w:val="false"/><w:b w:val="false"/><w:bCs w:val="false"/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr>'

function New-SimplePara([string]$text) {
    return '<w:p>' + $pPr36 + '<w:r>' + $rPr36 + '<w:t>' + $text + '</w:t></w:r></w:p>'
}

# ---------------------------------------------------------------------
# 1) Strip the __DdeLink__ bookmark around the first "Entities:" line.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$bookmarkFreeXml = New-PkgXml ('<w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:rFonts w:ascii="Ubuntu" w:hAnsi="Ubuntu"/><w:b/><w:b/><w:bCs/><w:sz w:val="52"/><w:szCs w:val="52"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Ubuntu" w:hAnsi="Ubuntu"/><w:b/><w:bCs/><w:sz w:val="52"/><w:szCs w:val="52"/></w:rPr><w:t>Entities:</w:t></w:r></w:p>')
$p1.Range.InsertXML($bookmarkFreeXml)

# ---------------------------------------------------------------------
# 2) Add the three new entity definitions after "Route(...)".
#    Route(...) is paragraph 5; reserve three empty paragraph slots right
#    after it (inheriting its sz=36 formatting) and then fill each one in
#    with InsertXML so we get exact control of the run/tab structure.
# ---------------------------------------------------------------------
$routeRange = $d.Paragraphs.Item(5).Range
$routeRange.Collapse(0)
$routeRange.InsertParagraphAfter()
$routeRange.Collapse(0)
$routeRange.InsertParagraphAfter()
$routeRange.Collapse(0)
$routeRange.InsertParagraphAfter()

$townXml    = New-PkgXml (New-SimplePara 'Town(name, garage)')
$garageXml  = New-PkgXml (New-SimplePara 'Garage(buses)')
$stageXml   = New-PkgXml ('<w:p>' + $pPr36 + '<w:r>' + $rPr36 + '<w:t>Stage</w:t></w:r><w:r>' + $rPr36 + '<w:t>(drivers)</w:t></w:r></w:p>')

$d.Paragraphs.Item(6).Range.InsertXML($townXml)
$d.Paragraphs.Item(7).Range.InsertXML($garageXml)
$d.Paragraphs.Item(8).Range.InsertXML($stageXml)

# ---------------------------------------------------------------------
# 3) Replace "Drivers drive Buses" / "Buses drive Routes" with the
#    expanded relationship list (5 paragraphs total).
# ---------------------------------------------------------------------
# Locate the two relationship paragraphs by their current text (formatting
# is identical to the rest of the body, sz=36 Ubuntu).
$driversPara = $null
$busesPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext -eq "Drivers drive Buses`r") { $driversPara = $i }
    if ($ptext -eq "Buses drive Routes`r") { $busesPara = $i }
}
Write-Host "driversPara index:" $driversPara " busesPara index:" $busesPara

# Reserve 3 additional empty paragraph slots right after "Buses drive
# Routes" (index $busesPara) so we end up with 5 total relationship
# paragraphs in the right place, then fill every slot via InsertXML.
$relRange = $d.Paragraphs.Item($busesPara).Range
$relRange.Collapse(0)
$relRange.InsertParagraphAfter()
$relRange.Collapse(0)
$relRange.InsertParagraphAfter()
$relRange.Collapse(0)
$relRange.InsertParagraphAfter()

$rel1 = New-PkgXml ('<w:p>' + $pPr36 + '<w:r>' + $rPr36 + '<w:t>Drivers allocated to Stages</w:t><w:tab/><w:t>many to one</w:t></w:r></w:p>')

$rel2 = New-PkgXml (
    '<w:p>' + $pPr36 +
    '<w:r>' + $rPr36 + '<w:t>Bus</w:t></w:r>' +
    '<w:r>' + $rPr36 + '<w:t>es</w:t></w:r>' +
    '<w:r>' + $rPr36 + '<w:t xml:space="preserve"> drive </w:t></w:r>' +
    '<w:r>' + $rPr36 + '<w:t>a</w:t></w:r>' +
    '<w:r>' + $rPr36 + '<w:t xml:space="preserve"> Rout</w:t></w:r>' +
    '<w:r>' + $rPr36 + '<w:t>e</w:t><w:tab/><w:tab/><w:tab/><w:t>many to one</w:t></w:r>' +
    '</w:p>'
)

$rel3 = New-PkgXml (
    '<w:p>' + $pPr36 +
    '<w:r>' + $rPr36 + '<w:t xml:space="preserve">Route has </w:t></w:r>' +
    '<w:r>' + $rPr36 + '<w:t>Stages</w:t><w:tab/><w:tab/><w:tab/><w:tab/><w:t>one to many</w:t></w:r>' +
    '</w:p>'
)

$rel4 = New-PkgXml (
    '<w:p>' + $pPr36 +
    '<w:r>' + $rPr36 + '<w:t>Town has a Garage</w:t><w:tab/><w:tab/><w:tab/></w:r>' +
    '<w:r>' + $rPr36 + '<w:t>one to one</w:t></w:r>' +
    '</w:p>'
)

$rel5 = New-PkgXml (
    '<w:p>' + $pPr36 +
    '<w:r>' + $rPr36 + '<w:t>Stage</w:t></w:r>' +
    '<w:r>' + $rPr36 + '<w:t xml:space="preserve"> through towns</w:t><w:tab/><w:tab/><w:tab/></w:r>' +
    '<w:r>' + $rPr36 + '<w:t>one to many</w:t></w:r>' +
    '</w:p>'
)

$d.Paragraphs.Item($driversPara).Range.InsertXML($rel1)
$d.Paragraphs.Item($busesPara).Range.InsertXML($rel2)
$d.Paragraphs.Item($busesPara + 1).Range.InsertXML($rel3)
$d.Paragraphs.Item($busesPara + 2).Range.InsertXML($rel4)
$d.Paragraphs.Item($busesPara + 3).Range.InsertXML($rel5)

Write-Host "Final paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host $i ":" ('[' + $d.Paragraphs.Item($i).Range.Text + ']')
}
